$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Walker Kessler", "C", "Utah Jazz")
)

$row = 2
foreach ($p in $players) {
    $ws.Cells.Item($row, 1).Value = $p[0]
    $ws.Cells.Item($row, 2).Value = $p[1]
    $ws.Cells.Item($row, 3).Value = $p[2]
    $row = $row + 1
}
